# Appends 22 new benchmark-result rows (575-596) to Sheet1 and touches up
# the previous last row (574): its E/F 'training/testing duration' cells
# lose the integer number format (now on the new last row, 596 instead) and
# its L/M timestamps get a tiny precision nudge.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 574 edits: drop the '0' integer number format from E574/F574 (it now
# belongs to the new last row, 596) and nudge the L574/M574 timestamps.
$ws.Cells.Item(574,5).ClearFormats()
$ws.Cells.Item(574,6).ClearFormats()
$ws.Cells.Item(574,12).Value = [double]"45929.47446090278"
$ws.Cells.Item(574,13).Value = [double]"45929.47446027778"

# Row 575
$ws.Cells.Item(575,1).Value = "Fucntionality_test_MUTAG_with_feature-KNN"
$ws.Cells.Item(575,2).Value = "MUTAG"
$ws.Cells.Item(575,3).Value = "feature-KNN"
$ws.Cells.Item(575,4).Value = [double]"0.2"
$ws.Cells.Item(575,5).Value = [double]"3.266898148148148e-07"
$ws.Cells.Item(575,6).Value = [double]"2.87962962962963e-08"
$ws.Cells.Item(575,7).Value = [double]"0.7368421052631579"
$ws.Cells.Item(575,8).Value = [double]"0.7172619047619048"
$ws.Cells.Item(575,9).Value = [double]"0.7403846153846154"
$ws.Cells.Item(575,10).Value = [double]"0.7130681818181819"
$ws.Cells.Item(575,11).Value = [double]"0.2471590909090909"
$ws.Cells.Item(575,12).Value = [double]"45934.69965002315"
$ws.Cells.Item(575,13).Value = [double]"45934.6996496875"
$ws.Cells.Item(575,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(575,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(575,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(575,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 576
$ws.Cells.Item(576,1).Value = "Fucntionality_test_MUTAG_with_SVC_Hybrid-Prototype-GED_poly"
$ws.Cells.Item(576,2).Value = "MUTAG"
$ws.Cells.Item(576,3).Value = "SVC_Hybrid-Prototype-GED_poly"
$ws.Cells.Item(576,4).Value = [double]"0.2"
$ws.Cells.Item(576,5).Value = [double]"3.484722222222222e-07"
$ws.Cells.Item(576,6).Value = [double]"2.459490740740741e-08"
$ws.Cells.Item(576,7).Value = [double]"0.8157894736842105"
$ws.Cells.Item(576,8).Value = [double]"0.7913725490196079"
$ws.Cells.Item(576,9).Value = [double]"0.7980769230769231"
$ws.Cells.Item(576,10).Value = [double]"0.7861538461538462"
$ws.Cells.Item(576,11).Value = [double]"0.1323076923076923"
$ws.Cells.Item(576,12).Value = [double]"45934.73216295139"
$ws.Cells.Item(576,13).Value = [double]"45934.73216259259"
$ws.Cells.Item(576,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(576,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(576,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(576,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 577
$ws.Cells.Item(577,1).Value = "Fucntionality_test_PTC_FR_with_SVC_Hybrid-Prototype-GED_poly"
$ws.Cells.Item(577,2).Value = "PTC_FR"
$ws.Cells.Item(577,3).Value = "SVC_Hybrid-Prototype-GED_poly"
$ws.Cells.Item(577,4).Value = [double]"0.2"
$ws.Cells.Item(577,5).Value = [double]"1.0690625e-06"
$ws.Cells.Item(577,6).Value = [double]"4.201388888888889e-08"
$ws.Cells.Item(577,7).Value = [double]"0.647887323943662"
$ws.Cells.Item(577,8).Value = [double]"0.3931623931623932"
$ws.Cells.Item(577,9).Value = [double]"0.323943661971831"
$ws.Cells.Item(577,10).Value = [double]"0.5"
$ws.Cells.Item(577,11).Value = [double]"0.5434782608695652"
$ws.Cells.Item(577,12).Value = [double]"45934.74649184028"
$ws.Cells.Item(577,13).Value = [double]"45934.74649076389"
$ws.Cells.Item(577,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(577,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(577,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(577,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 578
$ws.Cells.Item(578,1).Value = "Fucntionality_test_MUTAG_with_SVC_VertexHistogram_rbf"
$ws.Cells.Item(578,2).Value = "MUTAG"
$ws.Cells.Item(578,3).Value = "SVC_VertexHistogram_rbf"
$ws.Cells.Item(578,4).Value = [double]"0.2"
$ws.Cells.Item(578,5).Value = [double]"5.610362268518519e-05"
$ws.Cells.Item(578,6).Value = [double]"3.797453703703704e-08"
$ws.Cells.Item(578,7).Value = [double]"0.7894736842105263"
$ws.Cells.Item(578,8).Value = [double]"0.7088122605363985"
$ws.Cells.Item(578,9).Value = [double]"0.6902356902356902"
$ws.Cells.Item(578,10).Value = [double]"0.76036866359447"
$ws.Cells.Item(578,11).Value = [double]"0.1382488479262673"
$ws.Cells.Item(578,12).Value = [double]"45934.77812075231"
$ws.Cells.Item(578,13).Value = [double]"45934.77806451389"
$ws.Cells.Item(578,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(578,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(578,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(578,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 579
$ws.Cells.Item(579,1).Value = "Fucntionality_test_MUTAG_with_SVC_EdgeHistogram_rbf"
$ws.Cells.Item(579,2).Value = "MUTAG"
$ws.Cells.Item(579,3).Value = "SVC_EdgeHistogram_rbf"
$ws.Cells.Item(579,4).Value = [double]"0.2"
$ws.Cells.Item(579,5).Value = [double]"1.131828703703704e-07"
$ws.Cells.Item(579,6).Value = [double]"2.368055555555556e-08"
$ws.Cells.Item(579,7).Value = [double]"0.7105263157894737"
$ws.Cells.Item(579,8).Value = [double]"0.672156862745098"
$ws.Cells.Item(579,9).Value = [double]"0.7337164750957854"
$ws.Cells.Item(579,10).Value = [double]"0.6732954545454546"
$ws.Cells.Item(579,11).Value = [double]"0.1619318181818182"
$ws.Cells.Item(579,12).Value = [double]"45934.78180104167"
$ws.Cells.Item(579,13).Value = [double]"45934.78180091435"
$ws.Cells.Item(579,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(579,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(579,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(579,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 580
$ws.Cells.Item(580,1).Value = "Fucntionality_test_MUTAG_with_SVC_EdgeHistogram_rbf"
$ws.Cells.Item(580,2).Value = "MUTAG"
$ws.Cells.Item(580,3).Value = "SVC_EdgeHistogram_rbf"
$ws.Cells.Item(580,4).Value = [double]"0.2"
$ws.Cells.Item(580,5).Value = [double]"1.131828703703704e-07"
$ws.Cells.Item(580,6).Value = [double]"2.368055555555556e-08"
$ws.Cells.Item(580,7).Value = [double]"0.7105263157894737"
$ws.Cells.Item(580,8).Value = [double]"0.672156862745098"
$ws.Cells.Item(580,9).Value = [double]"0.7337164750957854"
$ws.Cells.Item(580,10).Value = [double]"0.6732954545454546"
$ws.Cells.Item(580,11).Value = [double]"0.6732954545454546"
$ws.Cells.Item(580,12).Value = [double]"45934.78180104167"
$ws.Cells.Item(580,13).Value = [double]"45934.78180091435"
$ws.Cells.Item(580,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(580,15).Value = "Hyperparameter Tuning (grid)"
$ws.Cells.Item(580,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(580,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 581
$ws.Cells.Item(581,1).Value = "Fucntionality_test_MUTAG_with_SVC_CombinedHistogram_rbf"
$ws.Cells.Item(581,2).Value = "MUTAG"
$ws.Cells.Item(581,3).Value = "SVC_CombinedHistogram_rbf"
$ws.Cells.Item(581,4).Value = [double]"0.2"
$ws.Cells.Item(581,5).Value = [double]"1.427777777777778e-07"
$ws.Cells.Item(581,6).Value = [double]"2.83912037037037e-08"
$ws.Cells.Item(581,7).Value = [double]"0.8157894736842105"
$ws.Cells.Item(581,8).Value = [double]"0.8048422597212033"
$ws.Cells.Item(581,9).Value = [double]"0.823076923076923"
$ws.Cells.Item(581,10).Value = [double]"0.7982954545454546"
$ws.Cells.Item(581,11).Value = [double]"0.1789772727272727"
$ws.Cells.Item(581,12).Value = [double]"45934.84572197917"
$ws.Cells.Item(581,13).Value = [double]"45934.8457218287"
$ws.Cells.Item(581,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(581,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(581,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(581,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 582
$ws.Cells.Item(582,1).Value = "Fucntionality_test_MUTAG_with_SVC_CombinedHistogram_rbf"
$ws.Cells.Item(582,2).Value = "MUTAG"
$ws.Cells.Item(582,3).Value = "SVC_CombinedHistogram_rbf"
$ws.Cells.Item(582,4).Value = [double]"0.2"
$ws.Cells.Item(582,5).Value = [double]"1.427777777777778e-07"
$ws.Cells.Item(582,6).Value = [double]"2.83912037037037e-08"
$ws.Cells.Item(582,7).Value = [double]"0.7631578947368421"
$ws.Cells.Item(582,8).Value = [double]"0.7414965986394557"
$ws.Cells.Item(582,9).Value = [double]"0.7794612794612794"
$ws.Cells.Item(582,10).Value = [double]"0.7357954545454546"
$ws.Cells.Item(582,11).Value = [double]"0.7357954545454546"
$ws.Cells.Item(582,12).Value = [double]"45934.84572197917"
$ws.Cells.Item(582,13).Value = [double]"45934.8457218287"
$ws.Cells.Item(582,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(582,15).Value = "Hyperparameter Tuning (grid)"
$ws.Cells.Item(582,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(582,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 583
$ws.Cells.Item(583,1).Value = "Fucntionality_test_MUTAG_with_SVC_CombinedHistogram_rbf"
$ws.Cells.Item(583,2).Value = "MUTAG"
$ws.Cells.Item(583,3).Value = "SVC_CombinedHistogram_rbf"
$ws.Cells.Item(583,4).Value = [double]"0.2"
$ws.Cells.Item(583,5).Value = [double]"1.173148148148148e-07"
$ws.Cells.Item(583,6).Value = [double]"2.293981481481481e-08"
$ws.Cells.Item(583,7).Value = [double]"0.7894736842105263"
$ws.Cells.Item(583,8).Value = [double]"0.7285714285714285"
$ws.Cells.Item(583,9).Value = [double]"0.775"
$ws.Cells.Item(583,10).Value = [double]"0.7115384615384616"
$ws.Cells.Item(583,11).Value = [double]"0.1378205128205128"
$ws.Cells.Item(583,12).Value = [double]"45934.85872717592"
$ws.Cells.Item(583,13).Value = [double]"45934.85872704861"
$ws.Cells.Item(583,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(583,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(583,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(583,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 584
$ws.Cells.Item(584,1).Value = "Fucntionality_test_MUTAG_with_SVC_CombinedHistogram_rbf"
$ws.Cells.Item(584,2).Value = "MUTAG"
$ws.Cells.Item(584,3).Value = "SVC_CombinedHistogram_rbf"
$ws.Cells.Item(584,4).Value = [double]"0.2"
$ws.Cells.Item(584,5).Value = [double]"1.173148148148148e-07"
$ws.Cells.Item(584,6).Value = [double]"2.293981481481481e-08"
$ws.Cells.Item(584,7).Value = [double]"0.7894736842105263"
$ws.Cells.Item(584,8).Value = [double]"0.7285714285714285"
$ws.Cells.Item(584,9).Value = [double]"0.775"
$ws.Cells.Item(584,10).Value = [double]"0.7115384615384616"
$ws.Cells.Item(584,11).Value = [double]"0.7115384615384616"
$ws.Cells.Item(584,12).Value = [double]"45934.85872717592"
$ws.Cells.Item(584,13).Value = [double]"45934.85872704861"
$ws.Cells.Item(584,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(584,15).Value = "Hyperparameter Tuning (grid)"
$ws.Cells.Item(584,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(584,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 585
$ws.Cells.Item(585,1).Value = "Fucntionality_test_MUTAG_with_SVC_CombinedHistogram_rbf"
$ws.Cells.Item(585,2).Value = "MUTAG"
$ws.Cells.Item(585,3).Value = "SVC_CombinedHistogram_rbf"
$ws.Cells.Item(585,4).Value = [double]"0.2"
$ws.Cells.Item(585,5).Value = [double]"0.0001553416898148148"
$ws.Cells.Item(585,6).Value = [double]"2.109953703703704e-07"
$ws.Cells.Item(585,7).Value = [double]"0.7368421052631579"
$ws.Cells.Item(585,8).Value = [double]"0.6801346801346801"
$ws.Cells.Item(585,9).Value = [double]"0.7834101382488479"
$ws.Cells.Item(585,10).Value = [double]"0.6782608695652175"
$ws.Cells.Item(585,11).Value = [double]"0.1318840579710145"
$ws.Cells.Item(585,12).Value = [double]"45934.86840295139"
$ws.Cells.Item(585,13).Value = [double]"45934.86824746527"
$ws.Cells.Item(585,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(585,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(585,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(585,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 586
$ws.Cells.Item(586,1).Value = "Fucntionality_test_MUTAG_with_feature-KNN"
$ws.Cells.Item(586,2).Value = "MUTAG"
$ws.Cells.Item(586,3).Value = "feature-KNN"
$ws.Cells.Item(586,4).Value = [double]"0.2"
$ws.Cells.Item(586,5).Value = [double]"4.058217592592592e-07"
$ws.Cells.Item(586,6).Value = [double]"4.969907407407407e-08"
$ws.Cells.Item(586,7).Value = [double]"0.7894736842105263"
$ws.Cells.Item(586,8).Value = [double]"0.7797101449275362"
$ws.Cells.Item(586,9).Value = [double]"0.7797101449275362"
$ws.Cells.Item(586,10).Value = [double]"0.7797101449275362"
$ws.Cells.Item(586,11).Value = [double]"0.136231884057971"
$ws.Cells.Item(586,12).Value = [double]"45934.87197950231"
$ws.Cells.Item(586,13).Value = [double]"45934.87197908565"
$ws.Cells.Item(586,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(586,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(586,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(586,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 587
$ws.Cells.Item(587,1).Value = "Fucntionality_test_MUTAG_with_feature-KNN"
$ws.Cells.Item(587,2).Value = "MUTAG"
$ws.Cells.Item(587,3).Value = "feature-KNN"
$ws.Cells.Item(587,4).Value = [double]"0.2"
$ws.Cells.Item(587,5).Value = [double]"4.006828703703703e-07"
$ws.Cells.Item(587,6).Value = [double]"4.663194444444444e-08"
$ws.Cells.Item(587,7).Value = [double]"0.7631578947368421"
$ws.Cells.Item(587,8).Value = [double]"0.7548387096774194"
$ws.Cells.Item(587,9).Value = [double]"0.7579710144927536"
$ws.Cells.Item(587,10).Value = [double]"0.7528409090909092"
$ws.Cells.Item(587,11).Value = [double]"0.1207386363636364"
$ws.Cells.Item(587,12).Value = [double]"45934.87640118055"
$ws.Cells.Item(587,13).Value = [double]"45934.87640077547"
$ws.Cells.Item(587,14).Value = "Dummy_Calculator"
$ws.Cells.Item(587,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(587,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(587,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 588
$ws.Cells.Item(588,1).Value = "Fucntionality_test_MUTAG_with_feature-KNN"
$ws.Cells.Item(588,2).Value = "MUTAG"
$ws.Cells.Item(588,3).Value = "feature-KNN"
$ws.Cells.Item(588,4).Value = [double]"0.2"
$ws.Cells.Item(588,5).Value = [double]"0.0003254243981481482"
$ws.Cells.Item(588,6).Value = [double]"8.603979166666666e-05"
$ws.Cells.Item(588,7).Value = [double]"0.9210526315789473"
$ws.Cells.Item(588,8).Value = [double]"0.9105882352941177"
$ws.Cells.Item(588,9).Value = [double]"0.9444444444444444"
$ws.Cells.Item(588,10).Value = [double]"0.8928571428571428"
$ws.Cells.Item(588,11).Value = [double]"0.0193452380952381"
$ws.Cells.Item(588,12).Value = [double]"45934.88721842592"
$ws.Cells.Item(588,13).Value = [double]"45934.88689296296"
$ws.Cells.Item(588,14).Value = "Dummy_Calculator"
$ws.Cells.Item(588,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(588,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(588,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 589
$ws.Cells.Item(589,1).Value = "Fucntionality_test_MUTAG_with_feature-KNN"
$ws.Cells.Item(589,2).Value = "MUTAG"
$ws.Cells.Item(589,3).Value = "feature-KNN"
$ws.Cells.Item(589,4).Value = [double]"0.2"
$ws.Cells.Item(589,5).Value = [double]"7.455092592592592e-07"
$ws.Cells.Item(589,6).Value = [double]"5.053240740740741e-08"
$ws.Cells.Item(589,7).Value = [double]"0.868421052631579"
$ws.Cells.Item(589,8).Value = [double]"0.8441345365053323"
$ws.Cells.Item(589,9).Value = [double]"0.9137931034482758"
$ws.Cells.Item(589,10).Value = [double]"0.8214285714285714"
$ws.Cells.Item(589,11).Value = [double]"0.06249999999999999"
$ws.Cells.Item(589,12).Value = [double]"45934.88872850694"
$ws.Cells.Item(589,13).Value = [double]"45934.88872775463"
$ws.Cells.Item(589,14).Value = "Dummy_Calculator"
$ws.Cells.Item(589,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(589,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(589,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 590
$ws.Cells.Item(590,1).Value = "Fucntionality_test_MUTAG_with_SVC_Hybrid-Prototype-GED_poly"
$ws.Cells.Item(590,2).Value = "MUTAG"
$ws.Cells.Item(590,3).Value = "SVC_Hybrid-Prototype-GED_poly"
$ws.Cells.Item(590,4).Value = [double]"0.2"
$ws.Cells.Item(590,5).Value = [double]"1.077604166666667e-06"
$ws.Cells.Item(590,6).Value = [double]"2.515277777777778e-07"
$ws.Cells.Item(590,7).Value = [double]"0.7105263157894737"
$ws.Cells.Item(590,8).Value = [double]"0.6933235509904623"
$ws.Cells.Item(590,9).Value = [double]"0.691304347826087"
$ws.Cells.Item(590,10).Value = [double]"0.6964285714285714"
$ws.Cells.Item(590,11).Value = [double]"0.2261904761904762"
$ws.Cells.Item(590,12).Value = [double]"45934.98598952546"
$ws.Cells.Item(590,13).Value = [double]"45934.98598822916"
$ws.Cells.Item(590,14).Value = "Dummy_Calculator"
$ws.Cells.Item(590,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(590,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(590,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 591
$ws.Cells.Item(591,1).Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(591,2).Value = "MUTAG"
$ws.Cells.Item(591,3).Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(591,4).Value = [double]"0.2"
$ws.Cells.Item(591,5).Value = [double]"4.554085648148148e-06"
$ws.Cells.Item(591,6).Value = [double]"2.150150462962963e-06"
$ws.Cells.Item(591,7).Value = [double]"0.8421052631578947"
$ws.Cells.Item(591,8).Value = [double]"0.7964285714285715"
$ws.Cells.Item(591,9).Value = [double]"0.7964285714285715"
$ws.Cells.Item(591,10).Value = [double]"0.7964285714285715"
$ws.Cells.Item(591,11).Value = [double]"0.1321428571428571"
$ws.Cells.Item(591,12).Value = [double]"45935.12126365741"
$ws.Cells.Item(591,13).Value = [double]"45935.12125907408"
$ws.Cells.Item(591,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(591,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(591,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(591,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 592
$ws.Cells.Item(592,1).Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(592,2).Value = "MUTAG"
$ws.Cells.Item(592,3).Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(592,4).Value = [double]"0.2"
$ws.Cells.Item(592,5).Value = [double]"4.452719907407407e-06"
$ws.Cells.Item(592,6).Value = [double]"3.83212962962963e-06"
$ws.Cells.Item(592,7).Value = [double]"0.7105263157894737"
$ws.Cells.Item(592,8).Value = [double]"0.5832502492522433"
$ws.Cells.Item(592,9).Value = [double]"0.8428571428571429"
$ws.Cells.Item(592,10).Value = [double]"0.6071428571428571"
$ws.Cells.Item(592,11).Value = [double]"0.1607142857142857"
$ws.Cells.Item(592,12).Value = [double]"45935.12517533565"
$ws.Cells.Item(592,13).Value = [double]"45935.12517086806"
$ws.Cells.Item(592,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(592,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(592,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(592,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 593
$ws.Cells.Item(593,1).Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(593,2).Value = "MUTAG"
$ws.Cells.Item(593,3).Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(593,4).Value = [double]"0.2"
$ws.Cells.Item(593,5).Value = [double]"1.090405092592593e-06"
$ws.Cells.Item(593,6).Value = [double]"3.418518518518519e-07"
$ws.Cells.Item(593,7).Value = [double]"0.7105263157894737"
$ws.Cells.Item(593,8).Value = [double]"0.672156862745098"
$ws.Cells.Item(593,9).Value = [double]"0.6692307692307693"
$ws.Cells.Item(593,10).Value = [double]"0.6762820512820513"
$ws.Cells.Item(593,11).Value = [double]"0.1858974358974359"
$ws.Cells.Item(593,12).Value = [double]"45935.12592638889"
$ws.Cells.Item(593,13).Value = [double]"45935.12592527777"
$ws.Cells.Item(593,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(593,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(593,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(593,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 594
$ws.Cells.Item(594,1).Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(594,2).Value = "MUTAG"
$ws.Cells.Item(594,3).Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(594,4).Value = [double]"0.2"
$ws.Cells.Item(594,5).Value = [double]"7.539930555555555e-07"
$ws.Cells.Item(594,6).Value = [double]"3.489236111111111e-07"
$ws.Cells.Item(594,7).Value = [double]"0.7368421052631579"
$ws.Cells.Item(594,8).Value = [double]"0.7076923076923077"
$ws.Cells.Item(594,9).Value = [double]"0.7023809523809523"
$ws.Cells.Item(594,10).Value = [double]"0.717948717948718"
$ws.Cells.Item(594,11).Value = [double]"0.201923076923077"
$ws.Cells.Item(594,12).Value = [double]"45935.12736097222"
$ws.Cells.Item(594,13).Value = [double]"45935.12736019676"
$ws.Cells.Item(594,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(594,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(594,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(594,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 595
$ws.Cells.Item(595,1).Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(595,2).Value = "MUTAG"
$ws.Cells.Item(595,3).Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(595,4).Value = [double]"0.2"
$ws.Cells.Item(595,5).Value = [double]"9.521296296296296e-07"
$ws.Cells.Item(595,6).Value = [double]"3.660300925925926e-07"
$ws.Cells.Item(595,7).Value = [double]"0.7631578947368421"
$ws.Cells.Item(595,8).Value = [double]"0.7414965986394557"
$ws.Cells.Item(595,9).Value = [double]"0.7461538461538462"
$ws.Cells.Item(595,10).Value = [double]"0.7380952380952381"
$ws.Cells.Item(595,11).Value = [double]"0.1220238095238095"
$ws.Cells.Item(595,12).Value = [double]"45935.70941407407"
$ws.Cells.Item(595,13).Value = [double]"45935.70941309028"
$ws.Cells.Item(595,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(595,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(595,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(595,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 596
$ws.Cells.Item(596,1).Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(596,2).Value = "MUTAG"
$ws.Cells.Item(596,3).Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Cells.Item(596,4).Value = [double]"0.2"
$ws.Cells.Item(596,5).Value = [double]"8.128356481481482e-07"
$ws.Cells.Item(596,6).Value = [double]"3.821990740740741e-07"
$ws.Cells.Item(596,7).Value = [double]"0.6842105263157895"
$ws.Cells.Item(596,8).Value = [double]"0.6761363636363636"
$ws.Cells.Item(596,9).Value = [double]"0.6750700280112045"
$ws.Cells.Item(596,10).Value = [double]"0.681159420289855"
$ws.Cells.Item(596,11).Value = [double]"0.3130434782608695"
$ws.Cells.Item(596,12).Value = [double]"45935.70991768294"
$ws.Cells.Item(596,13).Value = [double]"45935.70991684672"
$ws.Cells.Item(596,14).Value = "GEDLIB_Calculator"
$ws.Cells.Item(596,15).Value = "Simple Train-Test Split"
$ws.Cells.Item(596,12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(596,13).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(596,5).NumberFormat = "0"
$ws.Cells.Item(596,6).NumberFormat = "0"

